# Apply the commit: "chore: simulator full-month coverage, persist logs, fix employees"
#
# Changes:
#  - Weekly Timesheet sheet: fix three client/person names, bump hours from 8 -> 9
#    for each of the five daily rows (and their rate/total), and recompute the
#    SUBTOTAL / HOURLY SUBTOTAL / ADMIN SUBTOTAL / GRAND TOTAL rows.
#  - Jason Schema sheet: mirror the same name fixes, hour bump, and update the
#    Employee ID for every row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Weekly Timesheet"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Fix mis-spelled / corrected client names
$ws1.Range("B3").Value = "Evans"
$ws1.Range("B5").Value = "Hewett"
$ws1.Range("B6").Value = "Howard"

# Daily rows 2-6: Hours 8 -> 9, Rate 0 -> 140, Total 0 -> 1260
$dailyRows = 2, 3, 4, 5, 6
foreach ($r in $dailyRows) {
    $ws1.Cells.Item($r, 3).Value = 9
    $ws1.Cells.Item($r, 5).Value = 140
    $ws1.Cells.Item($r, 6).Value = 1260
}

# SUBTOTAL row (row 8): Hours 40 -> 45, Total 0 -> 6300, label updated
$ws1.Range("C8").Value = 45
$ws1.Range("F8").Value = 6300
$ws1.Range("D8").Value = "Reg: 45 / OT: 0"

# ADMIN SUBTOTAL (row 12) and GRAND TOTAL (row 13): Total 0 -> 6300
$ws1.Range("F12").Value = 6300
$ws1.Range("F13").Value = 6300

# ---------------------------------------------------------------------------
# Sheet 2: "Jason Schema"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Fix mis-spelled / corrected client names (column D)
$ws2.Range("D3").Value = "Evans"
$ws2.Range("D5").Value = "Hewett"
$ws2.Range("D6").Value = "Howard"

# Rows 2-6: Hours 8 -> 9, Rate 0 -> 140, Total 0 -> 1260, Employee ID update
$dataRows = 2, 3, 4, 5, 6
foreach ($r in $dataRows) {
    $ws2.Cells.Item($r, 2).Value = "emp_jp4mlvog"
    $ws2.Cells.Item($r, 5).Value = 9
    $ws2.Cells.Item($r, 6).Value = 140
    $ws2.Cells.Item($r, 7).Value = 1260
}
